$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1200053333333333
$ws.Range("H2").Value = 0.360016
$ws.Range("I2").Value = 0.01599335985210709
$ws.Range("J2").Value = 0.01792002013324337
$ws.Range("M2").Value = 162.7225033333333
$ws.Range("N2").Value = 488.16751
$ws.Range("O2").Value = 0.5231437953541009
$ws.Range("P2").Value = 0.5247717033381212
$ws.Range("Q2").Value = 19.52756825335111
$ws.Range("R2").Value = 175.74811428016
$ws.Range("S2").Value = 0.008366826973495205
$ws.Range("T2").Value = 0.009403919489175546
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1200053333333333
$ws.Range("H3").Value = 0.360016
$ws.Range("I3").Value = 0.01599335985210709
$ws.Range("J3").Value = 0.01792002013324337
$ws.Range("O3").Value = 0.0009322191998643353
$ws.Range("P3").Value = 0.0009351200601857102
$ws.Range("Q3").Value = 0.03479726647644445
$ws.Range("R3").Value = 0.313175398288
$ws.Range("S3").Value = 0.00001490931712447365
$ws.Range("T3").Value = 0.00001675737030552768
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1200053333333333
$ws.Range("H4").Value = 0.360016
$ws.Range("I4").Value = 0.01599335985210709
$ws.Range("J4").Value = 0.01792002013324337
$ws.Range("M4").Value = 61.580654
$ws.Range("N4").Value = 184.741962
$ws.Range("O4").Value = 0.1979783766474813
$ws.Range("P4").Value = 0.1985944416431287
$ws.Range("Q4").Value = 7.390006910154668
$ws.Range("R4").Value = 66.510062191392
$ws.Range("S4").Value = 0.003166339420659163
$ws.Range("T4").Value = 0.003558816392595091
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1200053333333333
$ws.Range("H5").Value = 0.360016
$ws.Range("I5").Value = 0.01599335985210709
$ws.Range("J5").Value = 0.01792002013324337
$ws.Range("M5").Value = 2.8947245
$ws.Range("N5").Value = 5.789449
$ws.Range("O5").Value = 0.009306378223129816
$ws.Range("P5").Value = 0.00622355841157717
$ws.Range("Q5").Value = 0.3473823785306667
$ws.Range("R5").Value = 2.084294271184
$ws.Range("S5").Value = 0.0001488402558423281
$ws.Range("T5").Value = 0.000111526292035879
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1200053333333333
$ws.Range("H6").Value = 0.360016
$ws.Range("I6").Value = 0.01599335985210709
$ws.Range("J6").Value = 0.01792002013324337
$ws.Range("M6").Value = 83.559527
$ws.Range("N6").Value = 250.678581
$ws.Range("O6").Value = 0.2686392305754237
$ws.Range("P6").Value = 0.2694751765469873
$ws.Range("Q6").Value = 10.02758889081067
$ws.Range("R6").Value = 90.248300017296
$ws.Range("S6").Value = 0.00429644388498592
$ws.Range("T6").Value = 0.004829000589131323
$ws.Range("I7").Value = 0.6614634286764723
$ws.Range("J7").Value = 0.7411474555000968
$ws.Range("M7").Value = 162.7225033333333
$ws.Range("N7").Value = 488.16751
$ws.Range("O7").Value = 0.5231437953541009
$ws.Range("P7").Value = 0.5247717033381212
$ws.Range("Q7").Value = 807.6334410041867
$ws.Range("R7").Value = 7268.700969037681
$ws.Range("S7").Value = 0.3460404885657464
$ws.Range("T7").Value = 0.3889332126475001
$ws.Range("I8").Value = 0.6614634286764723
$ws.Range("J8").Value = 0.7411474555000968
$ws.Range("O8").Value = 0.0009322191998643353
$ws.Range("P8").Value = 0.0009351200601857102
$ws.Range("S8").Value = 0.0006166289082203009
$ws.Range("T8").Value = 0.0006930618531937364
$ws.Range("I9").Value = 0.6614634286764723
$ws.Range("J9").Value = 0.7411474555000968
$ws.Range("M9").Value = 61.580654
$ws.Range("N9").Value = 184.741962
$ws.Range("O9").Value = 0.1979783766474813
$ws.Range("P9").Value = 0.1985944416431287
$ws.Range("Q9").Value = 305.640550449424
$ws.Range("R9").Value = 2750.764954044816
$ws.Range("S9").Value = 0.130955455821045
$ws.Range("T9").Value = 0.1471877651002673
$ws.Range("I10").Value = 0.6614634286764723
$ws.Range("J10").Value = 0.7411474555000968
$ws.Range("M10").Value = 2.8947245
$ws.Range("N10").Value = 5.789449
$ws.Range("O10").Value = 0.009306378223129816
$ws.Range("P10").Value = 0.00622355841157717
$ws.Range("Q10").Value = 14.367258742972
$ws.Range("R10").Value = 86.20355245783202
$ws.Range("S10").Value = 0.006155828848031504
$ws.Range("T10").Value = 0.004612574480896644
$ws.Range("I11").Value = 0.6614634286764723
$ws.Range("J11").Value = 0.7411474555000968
$ws.Range("M11").Value = 83.559527
$ws.Range("N11").Value = 250.678581
$ws.Range("O11").Value = 0.2686392305754237
$ws.Range("P11").Value = 0.2694751765469873
$ws.Range("Q11").Value = 414.727323739912
$ws.Range("R11").Value = 3732.545913659209
$ws.Range("S11").Value = 0.1776950265334292
$ws.Range("T11").Value = 0.199720841418239
$ws.Range("G12").Value = 2.420186
$ws.Range("H12").Value = 4.840372
$ws.Range("I12").Value = 0.3225432114714206
$ws.Range("J12").Value = 0.24093252436666
$ws.Range("M12").Value = 162.7225033333333
$ws.Range("N12").Value = 488.16751
$ws.Range("O12").Value = 0.5231437953541009
$ws.Range("P12").Value = 0.5247717033381212
$ws.Range("Q12").Value = 393.8187244522867
$ws.Range("R12").Value = 2362.91234671372
$ws.Range("S12").Value = 0.1687364798148593
$ws.Range("T12").Value = 0.1264345712014456
$ws.Range("G13").Value = 2.420186
$ws.Range("H13").Value = 4.840372
$ws.Range("I13").Value = 0.3225432114714206
$ws.Range("J13").Value = 0.24093252436666
$ws.Range("O13").Value = 0.0009322191998643353
$ws.Range("P13").Value = 0.0009351200601857102
$ws.Range("Q13").Value = 0.7017676200326667
$ws.Range("R13").Value = 4.210605720196001
$ws.Range("S13").Value = 0.0003006809745195608
$ws.Range("T13").Value = 0.0002253008366864462
$ws.Range("G14").Value = 2.420186
$ws.Range("H14").Value = 4.840372
$ws.Range("I14").Value = 0.3225432114714206
$ws.Range("J14").Value = 0.24093252436666
$ws.Range("M14").Value = 61.580654
$ws.Range("N14").Value = 184.741962
$ws.Range("O14").Value = 0.1979783766474813
$ws.Range("P14").Value = 0.1985944416431287
$ws.Range("Q14").Value = 149.036636681644
$ws.Range("R14").Value = 894.219820089864
$ws.Range("S14").Value = 0.06385658140577712
$ws.Range("T14").Value = 0.04784786015026634
$ws.Range("G15").Value = 2.420186
$ws.Range("H15").Value = 4.840372
$ws.Range("I15").Value = 0.3225432114714206
$ws.Range("J15").Value = 0.24093252436666
$ws.Range("M15").Value = 2.8947245
$ws.Range("N15").Value = 5.789449
$ws.Range("O15").Value = 0.009306378223129816
$ws.Range("P15").Value = 0.00622355841157717
$ws.Range("Q15").Value = 7.005771708757001
$ws.Range("R15").Value = 28.023086835028
$ws.Range("S15").Value = 0.003001709119255983
$ws.Range("T15").Value = 0.001499457638644648
$ws.Range("G16").Value = 2.420186
$ws.Range("H16").Value = 4.840372
$ws.Range("I16").Value = 0.3225432114714206
$ws.Range("J16").Value = 0.24093252436666
$ws.Range("M16").Value = 83.559527
$ws.Range("N16").Value = 250.678581
$ws.Range("O16").Value = 0.2686392305754237
$ws.Range("P16").Value = 0.2694751765469873
$ws.Range("Q16").Value = 202.229597412022
$ws.Range("R16").Value = 1213.377584472132
$ws.Range("S16").Value = 0.08664776015700858
$ws.Range("T16").Value = 0.06492533453961702
